# Practice IF, AND, OR, IFS functions
# Fill in the practice formulas on the "Data" sheet (columns F:O, rows 3-12)
# and the two goal-input cells (C15 Experience goal, C16 Total Salary goal).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Data")

# Goal inputs used by the L/M comparison formulas below
$ws.Range("C15").Value = 5
$ws.Range("C16").Value = 90000

for ($r = 3; $r -le 12; $r++) {
    # Total Salary (+)  ->  F : mirrors the Experience value
    $ws.Range("F$r").Formula = "=C$r"

    # Experience (=)    ->  G : Annual Salary + Bonus Max
    $ws.Range("G$r").Formula = "=D$r+E$r"

    # Bonus Rate (/)    ->  H : Bonus Max / Annual Salary, shown as a percentage
    $ws.Range("H$r").Formula = "=E$r/D$r"
    $ws.Range("H$r").NumberFormat = "0.0%"

    # Confirm Total Salary -> I : rebuild total salary from the bonus rate
    $ws.Range("I$r").Formula = "=H$r*D$r+D$r"

    # Does Total Salary = Confirmed Salary? -> J
    $ws.Range("J$r").Formula = "=I$r=G$r"

    # Is Bonus > Annual Salary? -> K
    $ws.Range("K$r").Formula = "=E$r>D$r"

    # Meets Experience -> L : experience goal (C15) <= person's experience
    $ws.Range("L$r").Formula = "=`$C`$15<=C$r"

    # Meets Salary -> M : total salary (G) >= salary goal (C16)
    $ws.Range("M$r").Formula = "=`$G$r>=C`$16"

    # Meets Both (1 or 0) -> N : AND of the two goals above, as 1/0
    $ws.Range("N$r").Formula = "=L$r*M$r"

    # Meets Both -> O : boolean version of N
    $ws.Range("O$r").Formula = "=N$r=1"
}
